$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.363018046572336
$ws.Range("C2").Value = 0.1319761253539582
$ws.Range("D2").Value = 0.1317583377187006
$ws.Range("E2").Value = 0.1135247653422233
$ws.Range("F2").Value = 1.491495209042583
$ws.Range("I2").Value = 0.8675160636340493
$ws.Range("J2").Value = 0.1293826194609364
$ws.Range("L2").Value = 0.3492390935779497
$ws.Range("N2").Value = 1.337447668177376
$ws.Range("O2").Value = 3.774698043901878

# Row 3
$ws.Range("B3").Value = 1.268414032269447
$ws.Range("C3").Value = 0.1183424876858794
$ws.Range("D3").Value = 0.1311258363583363
$ws.Range("E3").Value = 0.114054448640287
$ws.Range("F3").Value = 1.494010841709766
$ws.Range("I3").Value = 0.8749475156127176
$ws.Range("J3").Value = 0.1304495154349432
$ws.Range("L3").Value = 0.3418900136486656
$ws.Range("N3").Value = 1.346526100104356
$ws.Range("O3").Value = 3.784179701280891

# Row 4
$ws.Range("B4").Value = 1.210577653479788
$ws.Range("C4").Value = 0.1099183619519692
$ws.Range("D4").Value = 0.130774175444536
$ws.Range("E4").Value = 0.1144139746642381
$ws.Range("F4").Value = 1.496348933200629
$ws.Range("I4").Value = 0.8799707412763205
$ws.Range("J4").Value = 0.1311442677011723
$ws.Range("L4").Value = 0.3375047485299802
$ws.Range("N4").Value = 1.352589935462191
$ws.Range("O4").Value = 3.792141254159873

# Row 5
$ws.Range("B5").Value = 1.187073690941162
$ws.Range("C5").Value = 0.1064723040771867
$ws.Range("D5").Value = 0.1306401417731351
$ws.Range("E5").Value = 0.1145691258371553
$ws.Range("F5").Value = 1.497501328810429
$ws.Range("I5").Value = 0.8821334609038125
$ws.Range("J5").Value = 0.1314373718625061
$ws.Range("L5").Value = 0.3357498559432202
$ws.Range("N5").Value = 1.355184275717896
$ws.Range("O5").Value = 3.795923685846731

# Row 6
$ws.Range("B6").Value = 1.183174843045549
$ws.Range("C6").Value = 0.1058992991740837
$ws.Range("D6").Value = 0.1306184468004332
$ws.Range("E6").Value = 0.1145954109848084
$ws.Range("F6").Value = 1.497704741107462
$ws.Range("I6").Value = 0.8824995669082654
$ws.Range("J6").Value = 0.1314866450880405
$ws.Range("L6").Value = 0.3354604037076285
$ws.Range("N6").Value = 1.355622514239364
$ws.Range("O6").Value = 3.796584250873622

# Row 7
$ws.Range("B7").Value = 1.210260406002647
$ws.Range("C7").Value = 0.1098719402565109
$ws.Range("D7").Value = 0.1307723302240547
$ws.Range("E7").Value = 0.1144160320759733
$ws.Range("F7").Value = 1.496363666520011
$ws.Range("I7").Value = 0.8799994399841538
$ws.Range("J7").Value = 0.1311481801553729
$ws.Range("L7").Value = 0.3374809510700487
$ws.Range("N7").Value = 1.352624424287917
$ws.Range("O7").Value = 3.792190086920471

# Row 8
$ws.Range("B8").Value = 1.330347532172141
$ws.Range("C8").Value = 0.1272863508020521
$ws.Range("D8").Value = 0.1315326601373599
$ws.Range("E8").Value = 0.1137002915651024
$ws.Range("F8").Value = 1.492197954411871
$ws.Range("I8").Value = 0.869982879421368
$ws.Range("J8").Value = 0.1297422573776847
$ws.Range("L8").Value = 0.3466788552539128
$ws.Range("N8").Value = 1.340476377374117
$ws.Range("O8").Value = 3.777523182387455

# Row 9
$ws.Range("B9").Value = 1.56776383176873
$ws.Range("C9").Value = 0.1610092086410475
$ws.Range("D9").Value = 0.1333131096197278
$ws.Range("E9").Value = 0.112568146308357
$ws.Range("F9").Value = 1.490323190112704
$ws.Range("I9").Value = 0.8539940328816726
$ws.Range("J9").Value = 0.1272995380029407
$ws.Range("L9").Value = 0.3657177520956623
$ws.Range("N9").Value = 1.320532702456511
$ws.Range("O9").Value = 3.765745028524407

# Row 10
$ws.Range("B10").Value = 1.743298454634214
$ws.Range("C10").Value = 0.1855193228262806
$ws.Range("D10").Value = 0.1347955671645735
$ws.Range("E10").Value = 0.1119008765461427
$ws.Range("F10").Value = 1.492782960876013
$ws.Range("I10").Value = 0.8444766532774324
$ws.Range("J10").Value = 0.1256957449954541
$ws.Range("L10").Value = 0.3803091986813172
$ws.Range("N10").Value = 1.308236287807958
$ws.Range("O10").Value = 3.76745868650525

# Row 11
$ws.Range("B11").Value = 1.823380057555823
$ws.Range("C11").Value = 0.1966107164498396
$ws.Range("D11").Value = 0.1355074090039707
$ws.Range("E11").Value = 0.1116328437664151
$ws.Range("F11").Value = 1.494735361813724
$ws.Range("I11").Value = 0.84063162178699
$ws.Range("J11").Value = 0.1250074293808883
$ws.Range("L11").Value = 0.3870768737384509
$ws.Range("N11").Value = 1.303152315004695
$ws.Range("O11").Value = 3.770492679373092

# Row 12
$ws.Range("B12").Value = 1.853736419578581
$ws.Range("C12").Value = 0.2008021969646165
$ws.Range("D12").Value = 0.1357823120924664
$ws.Range("E12").Value = 0.1115364369966585
$ws.Range("F12").Value = 1.495594505083233
$ws.Range("I12").Value = 0.8392453196192378
$ws.Range("J12").Value = 0.1247527042713159
$ws.Range("L12").Value = 0.3896581472351954
$ws.Range("N12").Value = 1.3013003223043
$ws.Range("O12").Value = 3.771965926825033

# Row 13
$ws.Range("B13").Value = 1.847197276270379
$ws.Range("C13").Value = 0.1998998717568554
$ws.Range("D13").Value = 0.1357228698048587
$ws.Range("E13").Value = 0.1115569737166577
$ws.Range("F13").Value = 1.495404145143013
$ws.Range("I13").Value = 0.8395407824559626
$ws.Range("J13").Value = 0.1248073004843544
$ws.Range("L13").Value = 0.3891014038761824
$ws.Range("N13").Value = 1.301695928146692
$ws.Range("O13").Value = 3.771634209066178

# Row 14
$ws.Range("B14").Value = 1.825876879250302
$ws.Range("C14").Value = 0.1969557253344476
$ws.Range("D14").Value = 0.1355299185919989
$ws.Range("E14").Value = 0.1116248103698609
$ws.Range("F14").Value = 1.49480364305694
$ws.Range("I14").Value = 0.8405161720730945
$ws.Range("J14").Value = 0.1249863542899039
$ws.Range("L14").Value = 0.3872888669183681
$ws.Range("N14").Value = 1.302998484024982
$ws.Range("O14").Value = 3.770607383188434

# Row 15
$ws.Range("B15").Value = 1.812821527499239
$ws.Range("C15").Value = 0.1951512262594122
$ws.Range("D15").Value = 0.1354124250698518
$ws.Range("E15").Value = 0.1116670248859553
$ws.Range("F15").Value = 1.494451420092915
$ws.Range("I15").Value = 0.8411227088156323
$ws.Range("J15").Value = 0.1250968013065297
$ws.Range("L15").Value = 0.3861810396870453
$ws.Range("N15").Value = 1.303805866670196
$ws.Range("O15").Value = 3.770020665912483

# Row 16
$ws.Range("B16").Value = 1.738069401398093
$ws.Range("C16").Value = 0.184793282774109
$ws.Range("D16").Value = 0.1347497968015219
$ws.Range("E16").Value = 0.1119191064106957
$ws.Range("F16").Value = 1.492672134729503
$ws.Range("I16").Value = 0.8447376891870633
$ws.Range("J16").Value = 0.1257415575884817
$ws.Range("L16").Value = 0.3798695143441222
$ws.Range("N16").Value = 1.308578785692561
$ws.Range("O16").Value = 3.767305786278428

# Row 17
$ws.Range("B17").Value = 1.692268968716007
$ws.Range("C17").Value = 0.1784239288654703
$ws.Range("D17").Value = 0.1343528605155413
$ws.Range("E17").Value = 0.1120828348685716
$ws.Range("F17").Value = 1.491794038949266
$ws.Range("I17").Value = 0.8470794911240063
$ws.Range("J17").Value = 0.1261476562661077
$ws.Range("L17").Value = 0.3760307501452047
$ws.Range("N17").Value = 1.311637297256254
$ws.Range("O17").Value = 3.766217822271756

# Row 18
$ws.Range("B18").Value = 1.665947534154327
$ws.Range("C18").Value = 0.1747549615469666
$ws.Range("D18").Value = 0.1341280821782149
$ws.Range("E18").Value = 0.1121803505643051
$ws.Range("F18").Value = 1.491367430873652
$ws.Range("I18").Value = 0.8484720308091909
$ws.Range("J18").Value = 0.1263851173288657
$ws.Range("L18").Value = 0.3738350373447616
$ws.Range("N18").Value = 1.313444455934174
$ws.Range("O18").Value = 3.765804243674097

# Row 19
$ws.Range("B19").Value = 1.657039343830832
$ws.Range("C19").Value = 0.173511776403501
$ws.Range("D19").Value = 0.1340525835116892
$ws.Range("E19").Value = 0.1122139423993094
$ws.Range("F19").Value = 1.491236463563297
$ws.Range("I19").Value = 0.8489513503085888
$ws.Range("J19").Value = 0.1264661848561346
$ws.Range("L19").Value = 0.3730937150836979
$ws.Range("N19").Value = 1.314064573212242
$ws.Range("O19").Value = 3.765700652108251

# Row 20
$ws.Range("B20").Value = 1.697142264208026
$ws.Range("C20").Value = 0.1791025263851509
$ws.Range("D20").Value = 0.1343947501460576
$ws.Range("E20").Value = 0.1120650597732578
$ws.Range("F20").Value = 1.491879394784874
$ws.Range("I20").Value = 0.846825482730921
$ws.Range("J20").Value = 0.1261040244618208
$ws.Range("L20").Value = 0.3764381274035742
$ws.Range("N20").Value = 1.311306748087304
$ws.Range("O20").Value = 3.766311675950874

# Row 21
$ws.Range("B21").Value = 1.832138364955142
$ws.Range("C21").Value = 0.1978207271711767
$ws.Range("D21").Value = 0.1355864483501037
$ws.Range("E21").Value = 0.1116047470406034
$ws.Range("F21").Value = 1.494976773718705
$ws.Range("I21").Value = 0.8402277836174505
$ws.Range("J21").Value = 0.1249336011031961
$ws.Range("L21").Value = 0.3878207523001009
$ws.Range("N21").Value = 1.302613906010983
$ws.Range("O21").Value = 3.77090018311651

# Row 22
$ws.Range("B22").Value = 1.920547090840444
$ws.Range("C22").Value = 0.2100039969787986
$ws.Range("D22").Value = 0.1363964211121598
$ws.Range("E22").Value = 0.1113335758302654
$ws.Range("F22").Value = 1.497699425404818
$ws.Range("I22").Value = 0.8363222613351198
$ws.Range("J22").Value = 0.1242031941144361
$ws.Range("L22").Value = 0.3953677127063315
$ws.Range("N22").Value = 1.297359238563629
$ws.Range("O22").Value = 3.775789611324655

# Row 23
$ws.Range("B23").Value = 1.873345751389479
$ws.Range("C23").Value = 0.2035062101180358
$ws.Range("D23").Value = 0.1359612893701936
$ws.Range("E23").Value = 0.1114755951941131
$ws.Range("F23").Value = 1.496182409259319
$ws.Range("I23").Value = 0.8383695011217682
$ws.Range("J23").Value = 0.1245898688988447
$ws.Range("L23").Value = 0.391329959269811
$ws.Range("N23").Value = 1.30012475096364
$ws.Range("O23").Value = 3.773006993252352

# Row 24
$ws.Range("B24").Value = 1.694939016506225
$ws.Range("C24").Value = 0.1787957546759458
$ws.Range("D24").Value = 0.13437580117035
$ws.Range("E24").Value = 0.1120730853458216
$ws.Range("F24").Value = 1.491840561754856
$ws.Range("I24").Value = 0.8469401759981992
$ws.Range("J24").Value = 0.1261237379733249
$ws.Range("L24").Value = 0.3762539171085564
$ws.Range("N24").Value = 1.311456037411496
$ws.Range("O24").Value = 3.766268584633735

# Row 25
$ws.Range("B25").Value = 1.503337425576433
$ws.Range("C25").Value = 0.1519325806499694
$ws.Range("D25").Value = 0.1328006805200914
$ws.Range("E25").Value = 0.1128454624313573
$ws.Range("F25").Value = 1.490156521343934
$ws.Range("I25").Value = 0.8579280601583221
$ws.Range("J25").Value = 0.1279267868463849
$ws.Range("L25").Value = 0.3604607156110404
$ws.Range("N25").Value = 1.325513640109747
$ws.Range("O25").Value = 3.767111731833353

